# fix(publipostage): Add space before ":"
#
# The "statut_name" column contained "4: pas de résultats postés ni
# publiés" and needed a space before the colon: "4 : pas de résultats
# postés ni publiés". While re-generating the export, the source row
# order also shifted slightly (the "Role of Self-focused Attention in
# Depression" / NCT05464550 / RFASD / BEHAVIORAL record now comes before
# the NCT03304600 / tDCS-TOC / DEVICE record), and the Insula-TOP trial's
# intervention_type was corrected from DRUG to OTHER.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$correctedStatutName = "4 : pas de r" + [char]0x00E9 + "sultats post" + [char]0x00E9 + "s ni publi" + [char]0x00E9 + "s"

# Column B (statut_name) - fix spacing before colon for every data row
$ws.Range("B2:B6").Value = $correctedStatutName

# Column I (intervention_type) row 3: DRUG -> OTHER
$ws.Range("I3").Value = "OTHER"

# Rows 4 and 5: the NCT05464550 / "Role of Self-focused Attention in
# Depression" / RFASD / BEHAVIORAL record now occupies row 4 (it used to
# be row 5), and the NCT03304600 / tDCS record now occupies row 5.
$ws.Range("C4").Value = "NCT05464550"
$ws.Range("G4").Value = "Role of Self-focused Attention in Depression"
$ws.Range("H4").Value = "RFASD"
$ws.Range("I4").Value = "BEHAVIORAL"

$ws.Range("C5").Value = "NCT03304600"
$ws.Range("G5").Value = "Transcranial Direct Current Stimulation (tDCS) to Treat Patients With Severe and Resistant Obsessive Compulsive Disorder"
$ws.Range("H5").Value = "tDCS-TOC"
$ws.Range("I5").Value = "DEVICE"

$wb.Save()
